$wb = $excel.ActiveWorkbook

# --- Sheet "general": objValue/runtime/gap/Z-values updates ---
$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Cells.Item(3, 2).Value = 706.4655210009862
$wsGeneral.Cells.Item(4, 2).Value = 10307.1819999218
$wsGeneral.Cells.Item(5, 2).Value = [double]"5.662434058441323E-07"
$wsGeneral.Cells.Item(7, 2).Value = 0
$wsGeneral.Cells.Item(8, 2).Value = 0
$wsGeneral.Cells.Item(9, 2).Value = 257.8

# --- Sheet "x": update instance assignment values ---
$wsX = $wb.Worksheets.Item("x")
$wsX.Cells.Item(8, 2).Value = 12
$wsX.Cells.Item(10, 2).Value = 11

# --- Sheet "TBar": update target bar values ---
$wsTBar = $wb.Worksheets.Item("TBar")
$wsTBar.Cells.Item(6, 2).Value = 22.00864842470311
$wsTBar.Cells.Item(7, 2).Value = 21.98099133308861
$wsTBar.Cells.Item(8, 2).Value = 25.23527545392182
$wsTBar.Cells.Item(10, 2).Value = 10.28638906516763
$wsTBar.Cells.Item(11, 2).Value = 20
$wsTBar.Cells.Item(12, 2).Value = 18.45488985102369

# --- Sheet "Q": update C column values (rows 8-71) ---
$wsQ = $wb.Worksheets.Item("Q")
$wsQ.Cells.Item(8, 3).Value = 40.1
$wsQ.Cells.Item(9, 3).Value = 11.95
$wsQ.Cells.Item(10, 3).Value = 25.38
$wsQ.Cells.Item(11, 3).Value = 35.63
$wsQ.Cells.Item(12, 3).Value = 46.24
$wsQ.Cells.Item(13, 3).Value = 44.16
$wsQ.Cells.Item(14, 3).Value = 28.95
$wsQ.Cells.Item(15, 3).Value = 34.6
$wsQ.Cells.Item(16, 3).Value = 33.70999999999999
$wsQ.Cells.Item(17, 3).Value = 63.92
$wsQ.Cells.Item(18, 3).Value = 78.22
$wsQ.Cells.Item(19, 3).Value = 69.76000000000001
$wsQ.Cells.Item(20, 3).Value = 60.33
$wsQ.Cells.Item(21, 3).Value = 63.70999999999999
$wsQ.Cells.Item(22, 3).Value = 122.29
$wsQ.Cells.Item(23, 3).Value = 132.87
$wsQ.Cells.Item(24, 3).Value = 88.62
$wsQ.Cells.Item(25, 3).Value = 111.36
$wsQ.Cells.Item(26, 3).Value = 98.77
$wsQ.Cells.Item(27, 3).Value = 122.86
$wsQ.Cells.Item(28, 3).Value = 96.47
$wsQ.Cells.Item(29, 3).Value = 82.87
$wsQ.Cells.Item(30, 3).Value = 98.31999999999999
$wsQ.Cells.Item(31, 3).Value = 62.52000000000001
$wsQ.Cells.Item(32, 3).Value = 109.54
$wsQ.Cells.Item(33, 3).Value = 128.41
$wsQ.Cells.Item(34, 3).Value = 99.68000000000001
$wsQ.Cells.Item(35, 3).Value = 101.57
$wsQ.Cells.Item(36, 3).Value = 93.31999999999999
$wsQ.Cells.Item(37, 3).Value = 66.16
$wsQ.Cells.Item(38, 3).Value = 74.38
$wsQ.Cells.Item(39, 3).Value = 21.94
$wsQ.Cells.Item(40, 3).Value = 42.25
$wsQ.Cells.Item(41, 3).Value = 57.23999999999999
$wsQ.Cells.Item(42, 3).Value = 17.87
$wsQ.Cells.Item(43, 3).Value = 29.78
$wsQ.Cells.Item(44, 3).Value = 23.68
$wsQ.Cells.Item(45, 3).Value = 35.41
$wsQ.Cells.Item(46, 3).Value = 19.92
$wsQ.Cells.Item(47, 3).Value = 78.46000000000001
$wsQ.Cells.Item(48, 3).Value = 64.34999999999999
$wsQ.Cells.Item(49, 3).Value = 53.31
$wsQ.Cells.Item(50, 3).Value = 58.53
$wsQ.Cells.Item(51, 3).Value = 40.65000000000001
$wsQ.Cells.Item(52, 3).Value = 86.64
$wsQ.Cells.Item(53, 3).Value = 96.59999999999999
$wsQ.Cells.Item(54, 3).Value = 47.38
$wsQ.Cells.Item(55, 3).Value = 70.13
$wsQ.Cells.Item(56, 3).Value = 69.20999999999999
$wsQ.Cells.Item(57, 3).Value = 109.54
$wsQ.Cells.Item(58, 3).Value = 128.41
$wsQ.Cells.Item(59, 3).Value = 99.68000000000001
$wsQ.Cells.Item(60, 3).Value = 101.57
$wsQ.Cells.Item(61, 3).Value = 93.31999999999999
$wsQ.Cells.Item(62, 3).Value = 122.29
$wsQ.Cells.Item(63, 3).Value = 132.87
$wsQ.Cells.Item(64, 3).Value = 88.62
$wsQ.Cells.Item(65, 3).Value = 111.36
$wsQ.Cells.Item(66, 3).Value = 98.77
$wsQ.Cells.Item(67, 3).Value = 122.86
$wsQ.Cells.Item(68, 3).Value = 96.47
$wsQ.Cells.Item(69, 3).Value = 82.87
$wsQ.Cells.Item(70, 3).Value = 98.31999999999999
$wsQ.Cells.Item(71, 3).Value = 62.52000000000001

# --- Sheet "R": update C column values ---
$wsR = $wb.Worksheets.Item("R")
$wsR.Cells.Item(2, 3).Value = 9.539999999999992
$wsR.Cells.Item(3, 3).Value = 28.41
$wsR.Cells.Item(5, 3).Value = 1.569999999999993
$wsR.Cells.Item(7, 3).Value = 22.28999999999999
$wsR.Cells.Item(8, 3).Value = 32.87
$wsR.Cells.Item(10, 3).Value = 11.36
$wsR.Cells.Item(12, 3).Value = 22.86000000000001

# --- Sheet "y": remove row 2 (last instance row trimmed) ---
$wsY = $wb.Worksheets.Item("y")
$wsY.Rows.Item(2).Delete()

# --- Sheet "alpha": remove row 2 (last instance row trimmed) ---
$wsAlpha = $wb.Worksheets.Item("alpha")
$wsAlpha.Rows.Item(2).Delete()

